$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = 526
$ws.Range("B2").Value = 659
$ws.Range("B3").Value = 497
$ws.Range("B4").Value = 548
$ws.Range("B5").Value = 389
